$d = $word.ActiveDocument
$xml = $d.WordOpenXML
Write-Output "orig len=$($xml.Length)"
$newXml = $xml -replace "<w:docDefaults>.*?</w:docDefaults>", "REPLACED"
Write-Output "new len=$($newXml.Length)"
$idx = $newXml.IndexOf("REPLACED")
Write-Output "idx=$idx"
